# Apply cryptos list price/volume update (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.925.35'
$ws.Range("E2").Value = '  -0.06%  '

$ws.Range("D3").Value = '2.633.08'
$ws.Range("E3").Value = '  +4.07%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '110.65'
$ws.Range("E5").Value = '  +3.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '322.57'
$ws.Range("E6").Value = '  +0.79%  '

$ws.Range("E7").Value = '  -0.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.540'
$ws.Range("E9").Value = '  -1.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.56'
$ws.Range("E10").Value = '  -0.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.79'
$ws.Range("E11").Value = '  -1.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0810'
$ws.Range("E12").Value = '  +0.26%  '

$ws.Range("E13").Value = '  -0.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.21'
$ws.Range("E14").Value = '  +0.53%  '

$ws.Range("D15").Value = '3.037.26'
$ws.Range("E15").Value = '  +3.87%  '

$ws.Range("D16").Value = '2.635.05'
$ws.Range("E16").Value = '  +3.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.861'
$ws.Range("E17").Value = '  +1.42%  '

$ws.Range("D18").Value = '48.837.46'
$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  -1.21%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.67'
$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("B21").Value = 'ImmutableX'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.89'
$ws.Range("E21").Value = '  -1.87%  '

$ws.Range("D22").Value = '0.0₃0942'
$ws.Range("E22").Value = '  +0.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '271.21'
$ws.Range("E23").Value = '  -3.84%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.20'
$ws.Range("E24").Value = '  -2.82%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("E25").Value = '  +0.99%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.09'
$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.10'
$ws.Range("E28").Value = '  +3.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  +0.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.03'
$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("E31").Value = '  -4.13%  '

$ws.Range("E32").Value = '  -0.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.42'
$ws.Range("E33").Value = '  +1.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.21'
$ws.Range("E34").Value = '  -1.40%  '

$ws.Range("E35").Value = '  -0.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0795'
$ws.Range("E36").Value = '  +2.32%  '

$ws.Range("E37").Value = '  +6.56%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.03'
$ws.Range("E38").Value = '  +1.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.14'
$ws.Range("E39").Value = '  +6.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '124.62'
$ws.Range("E40").Value = '  +4.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.66'
$ws.Range("E41").Value = '  +2.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.111'
$ws.Range("E42").Value = '  -0.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.14'
$ws.Range("E43").Value = '  -3.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0313'
$ws.Range("E44").Value = '  +2.72%  '

$ws.Range("D45").Value = '2.065.86'
$ws.Range("E45").Value = '  +2.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.21'
$ws.Range("E46").Value = '  -0.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.11'
$ws.Range("E47").Value = '  +6.75%  '

$ws.Range("E48").Value = '  +4.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.93'
$ws.Range("E49").Value = '  -0.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '58.76'
$ws.Range("E50").Value = '  +3.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.17'
